$d = $word.ActiveDocument

# Locate the run of text that needs to be split into several runs
# (with w:proofErr gramStart/gramEnd markers inserted) and extended
# with the new sentence "Fazer a interface".
$target = $d.Content
$found = $target.Find.Execute(
    "persistência entender bem esse conceito da forma que o professor esta ensinando ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the paragraph text to edit"
}

# Re-seat a plain Range over the exact hit so InsertXML replaces the
# found text in place (a Range obtained straight off Find can behave
# like an insertion point instead of a replace target for InsertXML).
$r = $d.Range($target.Start, $target.End)

# The collapsed "_GoBack" bookmark sits immediately after this run; pull
# it out so it doesn't get dragged to the front of the replacement, then
# recreate it (still collapsed) right after the new, final run once the
# text is in place.
$d.Bookmarks("_GoBack").Delete()

$rPrXml = '<w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:b/><w:kern w:val="0"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:eastAsia="pt-BR"/><w14:ligatures w14:val="none"/></w:rPr>'

$bodyXml = (
    '<w:r>' + $rPrXml + '<w:t xml:space="preserve">persistência entender bem esse conceito da forma que o professor esta </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r>' + $rPrXml + '<w:t xml:space="preserve">ensinando </w:t></w:r>' +
    '<w:r>' + $rPrXml + '<w:t>.</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r>' + $rPrXml + '<w:t xml:space="preserve">  Fazer a interface</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>'
)

$packageXml = (
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body><w:p>' + $bodyXml + '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
)

$r.InsertXML($packageXml)
